$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D header text: "Prodi_Id" -> "Kode Prodi"
$ws.Range("D1").Value = "Kode Prodi"

# Remove the leftover placeholder cells in column C (rows 2-6), content + formatting
$ws.Range("C2:C6").Clear()

# Column A previously carried a centred cell style; drop back to default formatting
$ws.Columns.Item(1).ClearFormats()

# The unused Hyperlink cell style is no longer referenced anywhere - drop it
$wb.Styles.Item("Hyperlink").Delete()

# Re-size the columns to match the new template layout
$ws.Columns.Item(1).ColumnWidth = 12.67
$ws.Columns.Item(2).ColumnWidth = 15.67
$ws.Columns.Item(3).ColumnWidth = 15.67
$ws.Columns.Item(4).ColumnWidth = 10.5
$ws.Columns.Item(5).ColumnWidth = 10.67

# Match the saved selection state
$ws.Range("B2").Select() | Out-Null
